$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated power-flow results for the 380 kV case (vm_pu.xlsx, res_bus sheet).
# Maps each A1 cell reference (rows 2-25, columns B-N excluding empty H) to its new value.
$updates = @{
    "B2" = 1.02
    "C2" = 1.023853037233885
    "D2" = 1.046426690281516
    "E2" = 1.024418202208154
    "F2" = 1.051806270911759
    "I2" = 1.038597822941348
    "J2" = 1.029030929428512
    "K2" = 1.049192158254963
    "L2" = 1.027246762104584
    "M2" = 1.054556752670492
    "N2" = 1.013728238835013
    "B3" = 1.02
    "C3" = 1.024614579614752
    "D3" = 1.04697835063581
    "E3" = 1.025059803085877
    "F3" = 1.052501830847125
    "I3" = 1.038745142137015
    "J3" = 1.029432101398863
    "K3" = 1.049556134688157
    "L3" = 1.027695900305702
    "M3" = 1.055065333496253
    "N3" = 1.013862221798008
    "B4" = 1.02
    "C4" = 1.025108151980992
    "D4" = 1.047335907651674
    "E4" = 1.025476046816933
    "F4" = 1.052952894336875
    "I4" = 1.038839669537429
    "J4" = 1.029691812000185
    "K4" = 1.049791539015942
    "L4" = 1.027986920123492
    "M4" = 1.055394729579753
    "N4" = 1.01394892638787
    "B5" = 1.02
    "C5" = 1.025315840222298
    "D5" = 1.047486365194926
    "E5" = 1.025651293434718
    "F5" = 1.053142756048976
    "I5" = 1.038879217011232
    "J5" = 1.029801023048343
    "K5" = 1.049890474503355
    "L5" = 1.028109358504893
    "M5" = 1.055533280126733
    "N5" = 1.013985378584614
    "B6" = 1.02
    "C6" = 1.025350723101563
    "D6" = 1.047511635861204
    "E6" = 1.02568073316916
    "F6" = 1.053174648356761
    "I6" = 1.038885845932676
    "J6" = 1.029819361708663
    "K6" = 1.049907084492698
    "L6" = 1.028129921896475
    "M6" = 1.055556547566958
    "N6" = 1.013991499147817
    "B7" = 1.02
    "C7" = 1.02511092637566
    "D7" = 1.047337917522997
    "E7" = 1.025478387460559
    "F7" = 1.052955430359382
    "I7" = 1.038840198727541
    "J7" = 1.029693271172418
    "K7" = 1.049792861110604
    "L7" = 1.027988555785019
    "M7" = 1.055396580616655
    "N7" = 1.01394941345828
    "B8" = 1.02
    "C8" = 1.024110236501439
    "D8" = 1.046613001850945
    "E8" = 1.024634808350747
    "F8" = 1.052041132333341
    "I8" = 1.038647774891168
    "J8" = 1.029166480583135
    "K8" = 1.049315188243897
    "L8" = 1.027398467294819
    "M8" = 1.054728564585191
    "N8" = 1.013773516921723
    "B9" = 1.02
    "C9" = 1.022353133330482
    "D9" = 1.045340263471939
    "E9" = 1.023156713416005
    "F9" = 1.050437704694547
    "I9" = 1.038302623882868
    "J9" = 1.028239228635604
    "K9" = 1.048472670158234
    "L9" = 1.02636176483822
    "M9" = 1.053553893113842
    "N9" = 1.013463652449037
    "B10" = 1.02
    "C10" = 1.021186043425238
    "D10" = 1.044495025980832
    "E10" = 1.022177084468558
    "F10" = 1.049374056777309
    "I10" = 1.038068486911497
    "J10" = 1.027621826724518
    "K10" = 1.047910544894086
    "L10" = 1.025672808749539
    "M10" = 1.052772542231874
    "N10" = 1.013257164087596
    "B11" = 1.02
    "C11" = 1.020681728132174
    "D11" = 1.04412982728349
    "E11" = 1.021754287945052
    "F11" = 1.048914772973975
    "I11" = 1.037966155977153
    "J11" = 1.027354683555181
    "K11" = 1.047667051576662
    "L11" = 1.025375018784589
    "M11" = 1.052434650464642
    "N11" = 1.013167779358913
    "B12" = 1.02
    "C12" = 1.020494561586682
    "D12" = 1.04399429779027
    "E12" = 1.021597453312404
    "F12" = 1.048744369654748
    "I12" = 1.037928004101164
    "J12" = 1.027255485503811
    "K12" = 1.047576595439336
    "L12" = 1.025264487880038
    "M12" = 1.052309210189686
    "N12" = 1.013134582310867
    "B13" = 1.02
    "C13" = 1.02053470223637
    "D13" = 1.044023363790902
    "E13" = 1.021631085299746
    "F13" = 1.04878091287619
    "I13" = 1.037936194202958
    "J13" = 1.027276762405925
    "K13" = 1.047595999100502
    "L13" = 1.025288193419446
    "M13" = 1.052336114456872
    "N13" = 1.013141702982332
    "B14" = 1.02
    "C14" = 1.020666253639934
    "D14" = 1.044118621880936
    "E14" = 1.021741319632585
    "F14" = 1.048900683387812
    "I14" = 1.037963005215207
    "J14" = 1.027346483176693
    "K14" = 1.047659574677486
    "L14" = 1.025365880593123
    "M14" = 1.052424280138766
    "N14" = 1.013165035187511
    "B15" = 1.02
    "C15" = 1.020747327878782
    "D15" = 1.044177329707161
    "E15" = 1.021809266641781
    "F15" = 1.04897450388235
    "I15" = 1.037979505624712
    "J15" = 1.027389444569341
    "K15" = 1.0476987441611
    "L15" = 1.025413757076546
    "M15" = 1.052478610963018
    "N15" = 1.013179411529812
    "B16" = 1.02
    "C16" = 1.021219535302667
    "D16" = 1.044519279944407
    "E16" = 1.022205173531016
    "F16" = 1.049404565163262
    "I16" = 1.038075258374254
    "J16" = 1.027639560382713
    "K16" = 1.047926702993018
    "L16" = 1.025692583440455
    "M16" = 1.052794976394463
    "N16" = 1.013263096850104
    "B17" = 1.02
    "C17" = 1.021516019025545
    "D17" = 1.0447339905928
    "E17" = 1.022453888893904
    "F17" = 1.049674676380678
    "I17" = 1.038135068300194
    "J17" = 1.0277965048554
    "K17" = 1.048069672713389
    "L17" = 1.025867627520022
    "M17" = 1.052993542778119
    "N17" = 1.013315597779758
    "B18" = 1.02
    "C18" = 1.021689053400592
    "D18" = 1.044859304214769
    "E18" = 1.022599094261572
    "F18" = 1.049832351290326
    "I18" = 1.038169862886497
    "J18" = 1.027888066723198
    "K18" = 1.048153055805909
    "L18" = 1.025969779052358
    "M18" = 1.053109405353282
    "N18" = 1.013346223168047
    "B19" = 1.02
    "C19" = 1.021748070625381
    "D19" = 1.044902045836857
    "E19" = 1.022648628188121
    "F19" = 1.049886135252117
    "I19" = 1.038181711390738
    "J19" = 1.027919290102828
    "K19" = 1.048181485758318
    "L19" = 1.026004618736671
    "M19" = 1.053148918590405
    "N19" = 1.013356666040185
    "B20" = 1.02
    "C20" = 1.021484198724967
    "D20" = 1.044710946240401
    "E20" = 1.022427190223794
    "F20" = 1.049645683200759
    "I20" = 1.038128660722757
    "J20" = 1.027779664243013
    "K20" = 1.048054334306904
    "L20" = 1.025848841631455
    "M20" = 1.05297223410941
    "N20" = 1.013309964666308
    "B21" = 1.02
    "C21" = 1.020627510635886
    "D21" = 1.044090567374548
    "E21" = 1.021708852532555
    "F21" = 1.048865408551297
    "I21" = 1.037955113942843
    "J21" = 1.027325951285692
    "K21" = 1.047640853570099
    "L21" = 1.025343001393348
    "M21" = 1.052398315656683
    "N21" = 1.013158164308602
    "B22" = 1.02
    "C22" = 1.020089795491041
    "D22" = 1.043701214955712
    "E22" = 1.021258425897307
    "F22" = 1.048375949092359
    "I22" = 1.037845179288803
    "J22" = 1.027040862917145
    "K22" = 1.047380813864023
    "L22" = 1.025025432126066
    "M22" = 1.052037863402417
    "N22" = 1.013062747218541
    "B23" = 1.02
    "C23" = 1.020374760713552
    "D23" = 1.043907550549076
    "E23" = 1.021497089163446
    "F23" = 1.048635312879307
    "I23" = 1.037903535075108
    "J23" = 1.027191976230159
    "K23" = 1.047518671757223
    "L23" = 1.025193736297306
    "M23" = 1.052228908065213
    "N23" = 1.013113327009278
    "B24" = 1.02
    "C24" = 1.021498576639051
    "D24" = 1.044721358755143
    "E24" = 1.022439253790135
    "F24" = 1.049658783589178
    "I24" = 1.03813155631412
    "J24" = 1.027787273732917
    "K24" = 1.048061265099737
    "L24" = 1.025857330007626
    "M24" = 1.052981862449566
    "N24" = 1.013312510020193
    "B25" = 1.02
    "C25" = 1.022806635107043
    "D25" = 1.045668732277741
    "E25" = 1.023537828603464
    "F25" = 1.050851304475027
    "I25" = 1.038392568679809
    "J25" = 1.028478816322945
    "K25" = 1.048690566538802
    "L25" = 1.026629399714745
    "M25" = 1.053857271126124
    "N25" = 1.013543746460402
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

